# Update the "ASINs to Forecast" list with a new, shorter set of ASINs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New list of ASINs replacing the old one.
$asins = @(
    "B0C4BH2Z7D",
    "B0C4BKJCWV",
    "B0CDNCNBSP",
    "B0CDP24LW9",
    "B091ZG1ZC6",
    "B09SBV8M5F",
    "B0CXZK5JGC",
    "B08XY93JT3",
    "B07G5RD1XK",
    "B09MDTM5Z6",
    "B08FCY3BM2",
    "B0BHM59TQB",
    "B07TJX83W2",
    "B07QJ756H8",
    "B083GTQPXF",
    "B07TJWZGL9",
    "B084KGJDDM",
    "B099ZF7M85",
    "B0CKHM2ZQ6",
    "B07GBZL93X",
    "B08FWY81LM",
    "B07Z9YBT3T",
    "B091ZFTQQ1",
    "B07QBVWCQL",
    "B07GJKR7RX",
    "B08XY9146V",
    "B07VV87QVL"
)

# Clear out the previously used range (old data went down to row 45).
$ws.Range("A1:A45").ClearContents()

# Header stays the same.
$ws.Range("A1").Value = "ASIN"

# Write the new ASIN values into A2:A28.
for ($i = 0; $i -lt $asins.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $asins[$i]
}

# Match the saved selection/active cell from the edited workbook.
$ws.Range("C3").Select()
